$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the last used row of the sheet.
$lastRow = $ws.UsedRange.Rows.Count

# Column D holds "codeforiati:group-name" (header + values) and column E
# holds "codeforiati:group-code" (header + values). Swap the two columns,
# header row included, so D becomes group-code and E becomes group-name.
$dRange = $ws.Range("D1:D$lastRow")
$eRange = $ws.Range("E1:E$lastRow")

$dValues = $dRange.Value2
$eValues = $eRange.Value2

$dRange.Value = $eValues
$eRange.Value = $dValues
